$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 158, shifting existing rows 158-166 down to 159-167.
$ws.Rows.Item(158).Insert()

# Fill the new row 158 with data. Columns A,B,C,E,F,G,I,O,R are the same as the
# (old) row 158 which is now row 159; D,H,J,K,L,M,N,P,Q are new values.
$ws.Range("A158").Value = 10
$ws.Range("B158").Value = "Vega Modelo de Temuco"
$ws.Range("C158").Value = "La Araucanía"
$ws.Range("D158").Value = 44746
$ws.Range("D158").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E158").Value = 9
$ws.Range("F158").Value = 100112013
$ws.Range("G158").Value = "Alcachofa"
$ws.Range("H158").Value = "Española"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 85
$ws.Range("K158").Value = 22000
$ws.Range("L158").Value = 22000
$ws.Range("M158").Value = 22000
$ws.Range("N158").Value = "$/caja 30 unidades"
$ws.Range("O158").Value = "Provincia de Limarí"
$ws.Range("P158").Value = 733
$ws.Range("Q158").Value = 30
$ws.Range("R158").Value = "Hortaliza"
